# "Colocando header nos gráficos"
# Adds a header label to column A (row 1) of each data sheet that feeds a
# chart, fixes accented Portuguese text that had been stripped of
# diacritics, restyles the row-1 header cells (style "Normal" -> plain),
# and updates the "Custo Total" sheet's numbers/labels. Also drops the
# now-unused "Teto" row from the "Emissoes Totais" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share the same "Fonte/Tecnologia" table layout.
# ---------------------------------------------------------------------
$sourceSheets = @(1, 2, 3, 4)
foreach ($idx in $sourceSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # New header cell in A1, styled like the rest of the header row.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("A1").Style = "Header"

    # Row labels lose the bold/border header style and get accents fixed.
    $ws.Range("A2").Style = "Normal"
    $ws.Range("A2").Value = "Hidro"

    $ws.Range("A3").Style = "Normal"
    $ws.Range("A3").Value = "Gás Natural"

    $ws.Range("A4").Style = "Normal"
    $ws.Range("A4").Value = "Carvão"

    $ws.Range("A5").Style = "Normal"
    $ws.Range("A5").Value = "Nuclear"

    $ws.Range("A6").Style = "Normal"
    $ws.Range("A6").Value = "Óleos Comb"

    $ws.Range("A7").Style = "Normal"
    $ws.Range("A7").Value = "Biomassa"

    $ws.Range("A8").Style = "Normal"
    $ws.Range("A8").Value = "Eólica"

    $ws.Range("A9").Style = "Normal"
    $ws.Range("A9").Value = "Solar"

    $ws.Range("A10").Style = "Normal"
    $ws.Range("A10").Value = "Outros"

    $ws.Range("A11").Style = "Normal"
    $ws.Range("A11").Value = "Pot. Compl."

    $ws.Range("A12").Style = "Normal"
    $ws.Range("A12").Value = "GD"
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)" - add header, fix labels, drop
# the "Teto" row entirely.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Período"
$ws5.Range("A1").Style = "Header"

$ws5.Range("A2").Style = "Normal"
$ws5.Range("A2").Value = "P.Médio"

$ws5.Range("A3").Style = "Normal"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)" - add header, fix labels,
# update the expansion-cost values and the B1 column label.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("A1").Style = "Header"

$ws6.Range("B1").Value = "2015"

$ws6.Range("A2").Style = "Normal"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 169

$ws6.Range("A3").Style = "Normal"
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
